$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Column A (NC)
$ws.Cells.Item(2, 1).Value = 20330051920082
$ws.Cells.Item(3, 1).Value = 18330051920108
$ws.Cells.Item(4, 1).Value = 20330051920069

# Column B (Paterno)
$ws.Cells.Item(2, 2).Value = "GONZALEZ"
$ws.Cells.Item(3, 2).Value = "RAMIREZ"
$ws.Cells.Item(4, 2).Value = "VELAZQUEZ"

# Column C (Materno)
$ws.Cells.Item(2, 3).Value = "OFICIAL"
$ws.Cells.Item(3, 3).Value = "CORDOBA"
$ws.Cells.Item(4, 3).Value = "LOPEZ"

# Column D (Nombres)
$ws.Cells.Item(2, 4).Value = "SAID ANDRES"
$ws.Cells.Item(3, 4).Value = "FRANCISCO JAVIER"
$ws.Cells.Item(4, 4).Value = "AYLIN MELISSA"

# Column E (Nombre_Largo)
$ws.Cells.Item(2, 5).Value = "ÉTICA"
$ws.Cells.Item(3, 5).Value = "SUPERVISA EL CUMPLIMIENTO DE TAREAS Y PROCESOS PARA EVALUAR LA PRODUCTIVIDAD EN LA ORGANIZACIÓN"
$ws.Cells.Item(4, 5).Value = "ÉTICA"

# Column F (Grupo)
$ws.Cells.Item(2, 6).Value = "3AEV"
$ws.Cells.Item(3, 6).Value = "5ARHV"
$ws.Cells.Item(4, 6).Value = "3AEV"

# Column G (Reprobadas)
$ws.Cells.Item(2, 7).Value = 6
$ws.Cells.Item(3, 7).Value = 3
$ws.Cells.Item(4, 7).Value = 6
